$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9658119658119658
$ws.Range("C2").Value = 0.8154639175257732
$ws.Range("D2").Value = 0.8842929010620458
$ws.Range("E2").Value = 970

$ws.Range("B3").Value = 0.3035019455252918
$ws.Range("C3").Value = 0.7358490566037735
$ws.Range("D3").Value = 0.4297520661157024

$ws.Range("B4").Value = 0.8076208178438662
$ws.Range("C4").Value = 0.8076208178438662
$ws.Range("D4").Value = 0.8076208178438662
$ws.Range("E4").Value = 0.8076208178438662

$ws.Range("B5").Value = 0.6346569556686288
$ws.Range("C5").Value = 0.7756564870647733
$ws.Range("D5").Value = 0.6570224835888742
$ws.Range("E5").Value = 1076

$ws.Range("B6").Value = 0.9005658113971076
$ws.Range("C6").Value = 0.8076208178438662
$ws.Range("D6").Value = 0.8395147147197481
$ws.Range("E6").Value = 1076
